$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume change % (column E) values
$ws.Range("D2").Value = "'29.255.32"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'1.901.80"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'326.44"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4647"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.3918"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "'0.07896"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'0.9897"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("D12").Value = "'1.903.20"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "'7.076"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "'5.745"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'0.06998"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'88.29"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "'0.000009978"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "'17.12"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'29.271.79"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'2.142.95"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Value = "'2.101"
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").Value = "'156.14"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'19.40"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "'5.980"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").Value = "'118.78"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'1.885"
$ws.Range("E30").Value = "  -5.63%  "
$ws.Range("D31").Value = "'0.09359"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "'0.9015"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").Value = "'5.264"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").Value = "'1.325"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "'0.05784"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'0.02088"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").Value = "'7.710"
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("D41").Value = "'0.5705"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "'0.1785"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").Value = "'9.710"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").Value = "'11.91"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").Value = "'0.5360"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").Value = "'2.172"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").Value = "'0.07025"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "'1.854"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'2.573"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "'113.20"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").Value = "'1.052"
$ws.Range("E51").Value = "  -1.31%  "
